$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 126, shifting existing rows 126:197 down to 128:199
$ws.Range("A126:A127").EntireRow.Insert()

# Populate the two newly inserted rows (126 and 127) with their data
$ws.Range("A126").Value = 11
$ws.Range("B126").Value = "Vega Monumental Concepción"
$ws.Range("C126").Value = "Bíobío"
$ws.Range("D126").Value = 44960
$ws.Range("E126").Value = 8
$ws.Range("F126").Value = 100112044
$ws.Range("G126").Value = "Perejil"
$ws.Range("H126").Value = "Sin especificar"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 100
$ws.Range("K126").Value = 700
$ws.Range("L126").Value = 800
$ws.Range("M126").Value = 750
$ws.Range("N126").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O126").Value = "Región de Ñuble"
$ws.Range("P126").Value = 750
$ws.Range("Q126").Value = 1
$ws.Range("R126").Value = "Hortaliza"

$ws.Range("A127").Value = 11
$ws.Range("B127").Value = "Vega Monumental Concepción"
$ws.Range("C127").Value = "Bíobío"
$ws.Range("D127").Value = 44960
$ws.Range("E127").Value = 8
$ws.Range("F127").Value = 100112044
$ws.Range("G127").Value = "Perejil"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Segunda"
$ws.Range("J127").Value = 50
$ws.Range("K127").Value = 600
$ws.Range("L127").Value = 600
$ws.Range("M127").Value = 600
$ws.Range("N127").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O127").Value = "Región de Ñuble"
$ws.Range("P127").Value = 600
$ws.Range("Q127").Value = 1
$ws.Range("R127").Value = "Hortaliza"
